$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing "naja" -> "xxx" string (shared string index 9, cell E2)
$ws.Range("E2").Value = "xxx"

# New cells added in row 2
$ws.Range("G2").Value = "hmm"

# New cells added in row 3
$ws.Range("E3").Value = "bissifiel"
$ws.Range("F3").Value = "aufirohrum"

# New cells added in row 4
$ws.Range("E4").Value = "adfg"
$ws.Range("F4").Value = "arhf"

# New cells added in row 5
$ws.Range("F5").Value = "haerfg"
$ws.Range("G5").Value = "ashfrfdg"

# New cells added in row 6
$ws.Range("E6").Value = "asdfhpoh"
$ws.Range("F6").Value = "apfsogh"

# Update the view: scroll so column B is the top-left visible column,
# and move the selection/active cell to F4.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F4").Select()
